$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 21:47"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6501041
$ws.Range("C4").Value = 15466
$ws.Range("D4").Value = 3775915
$ws.Range("E4").Value = 2531321
$ws.Range("G4").Value = 271
$ws.Range("H4").Value = 193805

# Row 10 - Sudafrica
$ws.Range("B10").Value = 640441
$ws.Range("C10").Value = 1079
$ws.Range("D10").Value = 567729
$ws.Range("E10").Value = 57626
$ws.Range("G10").Value = 82
$ws.Range("H10").Value = 15086

# Row 17 - Francia
$ws.Range("E17").Value = 216534
$ws.Range("G17").Value = 38
$ws.Range("H17").Value = 30764

# Row 23 - Irak
$ws.Range("B23").Value = 269578
$ws.Range("C23").Value = 4894
$ws.Range("D23").Value = 206324
$ws.Range("E23").Value = 55597
$ws.Range("G23").Value = 68
$ws.Range("H23").Value = 7657

# Row 24 - Alemania
$ws.Range("B24").Value = 254929
$ws.Range("C24").Value = 1304
$ws.Range("E24").Value = 18520
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 9409

# Row 28 - Israel
$ws.Range("B28").Value = 137159
$ws.Range("C28").Value = 3184
$ws.Range("D28").Value = 106996
$ws.Range("E28").Value = 29123
$ws.Range("G28").Value = 14
$ws.Range("H28").Value = 1040

# Row 62 - Suiza
$ws.Range("E62").Value = 5119
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 2018

# Row 139 - Sudan del Sur
$ws.Range("B139").Value = 2552
$ws.Range("C139").Value = 7
$ws.Range("E139").Value = 1213
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 49

# Row 142 - Reunion
$ws.Range("B142").Value = 2346
$ws.Range("C142").Value = 69
$ws.Range("E142").Value = 1020

# Row 166 - Republica del Chad
$ws.Range("B166").Value = 1045
$ws.Range("C166").Value = 5
$ws.Range("D166").Value = 927
$ws.Range("E166").Value = 39
